$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the changed Price (column D) cells to remain text so values
# such as "232.18" are not reinterpreted as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values.
$ws.Range("D2").Value = "41.889.64"
$ws.Range("E2").Value = "  -4.58%  "
$ws.Range("D3").Value = "2.238.89"
$ws.Range("E3").Value = "  -4.77%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "232.18"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  -6.12%  "
$ws.Range("D7").Value = "69.16"
$ws.Range("E7").Value = "  -5.11%  "
$ws.Range("D9").Value = "0.552"
$ws.Range("E9").Value = "  -8.27%  "
$ws.Range("D10").Value = "0.0975"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").Value = "57.90"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").Value = "35.48"
$ws.Range("E12").Value = "  +6.66%  "
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("E14").Value = "  -7.62%  "
$ws.Range("D15").Value = "2.574.95"
$ws.Range("E15").Value = "  -4.78%  "
$ws.Range("D16").Value = "14.93"
$ws.Range("E16").Value = "  -8.84%  "
$ws.Range("D17").Value = "0.854"
$ws.Range("E17").Value = "  -5.91%  "
$ws.Range("D18").Value = "2.239.93"
$ws.Range("E18").Value = "  -4.78%  "
$ws.Range("D19").Value = "41.847.53"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  -6.57%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.22"
$ws.Range("E21").Value = "  -7.30%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "73.10"
$ws.Range("E22").Value = "  -6.02%  "
$ws.Range("D23").Value = "235.82"
$ws.Range("E23").Value = "  -7.91%  "
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").Value = "2.34"
$ws.Range("E27").Value = "  -6.24%  "
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("D30").Value = "168.34"
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("D31").Value = "20.58"
$ws.Range("E31").Value = "  -8.85%  "
$ws.Range("E32").Value = "  -7.62%  "
$ws.Range("E33").Value = "  -7.48%  "
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "0.0711"
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("D36").Value = "4.74"
$ws.Range("E36").Value = "  -7.32%  "
$ws.Range("D37").Value = "3.58"
$ws.Range("E37").Value = "  -5.62%  "
$ws.Range("D38").Value = "22.02"
$ws.Range("E38").Value = "  +15.96%  "
$ws.Range("E39").Value = "  -5.54%  "
$ws.Range("D40").Value = "5.99"
$ws.Range("E40").Value = "  -6.75%  "
$ws.Range("D41").Value = "0.0264"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "66.51"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "4.95"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("D44").Value = "9.01"
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -8.19%  "
$ws.Range("D46").Value = "0.189"
$ws.Range("E46").Value = "  -6.33%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -7.11%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -6.64%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "4.27"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("E51").Value = "  +3.44%  "
